$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 742
$ws.Range("F6").Value = 2315
$ws.Range("F7").Value = 53
$ws.Range("F8").Value = 1733
$ws.Range("F9").Value = 2919
$ws.Range("F10").Value = 166
$ws.Range("F11").Value = 4368
$ws.Range("F12").Value = 377
$ws.Range("F13").Value = 205
$ws.Range("F15").Value = 552
$ws.Range("F16").Value = 261
$ws.Range("F17").Value = 9
$ws.Range("F18").Value = 155
$ws.Range("F20").Value = 102
$ws.Range("F21").Value = 302
$ws.Range("F22").Value = 4370
$ws.Range("F24").Value = 3644
$ws.Range("F27").Value = 552
$ws.Range("F28").Value = 4369
$ws.Range("F29").Value = 85
$ws.Range("F30").Value = 528
$ws.Range("F31").Value = 540
$ws.Range("F32").Value = 494

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 26
$ws.Range("F5").Value = 29

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 9

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 9
$ws.Range("F8").Value = 742
$ws.Range("F9").Value = 2315
$ws.Range("F10").Value = 53
$ws.Range("F11").Value = 1733
$ws.Range("F13").Value = 2919
$ws.Range("F14").Value = 166
$ws.Range("F15").Value = 4368
$ws.Range("F16").Value = 377
$ws.Range("F17").Value = 205
$ws.Range("F19").Value = 552
$ws.Range("F20").Value = 261
$ws.Range("F21").Value = 9
$ws.Range("F22").Value = 155
$ws.Range("F23").Value = 26
$ws.Range("F25").Value = 102
$ws.Range("F26").Value = 302
$ws.Range("F27").Value = 4370
$ws.Range("F29").Value = 3644
$ws.Range("F32").Value = 553
$ws.Range("F33").Value = 4369
$ws.Range("F34").Value = 85
$ws.Range("F35").Value = 528
$ws.Range("F36").Value = 540
$ws.Range("F37").Value = 494
$ws.Range("F39").Value = 29
